$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2636.182
$ws.Cells.Item(86, 9).Value = 2471.4285
$ws.Cells.Item(86, 11).Value = 2471.4285
$ws.Cells.Item(86, 13).Value = -1348.4285
$ws.Cells.Item(89, 8).Value = 2636.182
$ws.Cells.Item(89, 9).Value = 2471.4285
$ws.Cells.Item(89, 11).Value = 12357.1425
$ws.Cells.Item(89, 13).Value = -6741.1425
$ws.Cells.Item(137, 8).Value = 13897414
$ws.Cells.Item(137, 10).Value = 16103.223
$ws.Cells.Item(137, 12).Value = 48309.669
$ws.Cells.Item(137, 14).Value = -53409.669

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 18183278
$ws.Cells.Item(2, 9).Value = 30303796
$ws.Cells.Item(2, 11).Value = 30303796
$ws.Cells.Item(2, 13).Value = -30303683
$ws.Cells.Item(45, 8).Value = 2439.389
$ws.Cells.Item(45, 9).Value = 2320.7334
$ws.Cells.Item(45, 11).Value = 2320.7334
$ws.Cells.Item(45, 13).Value = -1943.7334
$ws.Cells.Item(61, 9).Value = 4979.636
$ws.Cells.Item(61, 10).Value = 6939.125
$ws.Cells.Item(61, 11).Value = 4979.636
$ws.Cells.Item(61, 12).Value = 6939.125
$ws.Cells.Item(61, 13).Value = -4767.636
$ws.Cells.Item(61, 14).Value = -7363.125
$ws.Cells.Item(74, 8).Value = 2121.6
$ws.Cells.Item(74, 9).Value = 526.125
$ws.Cells.Item(74, 10).Value = 3001.862
$ws.Cells.Item(74, 11).Value = 526.125
$ws.Cells.Item(74, 12).Value = 3001.862
$ws.Cells.Item(74, 13).Value = 347.875
$ws.Cells.Item(74, 14).Value = -4749.862
$ws.Cells.Item(77, 8).Value = 2121.6
$ws.Cells.Item(77, 9).Value = 526.125
$ws.Cells.Item(77, 10).Value = 3001.862
$ws.Cells.Item(77, 11).Value = 2630.625
$ws.Cells.Item(77, 12).Value = 15009.31
$ws.Cells.Item(77, 13).Value = 1737.375
$ws.Cells.Item(77, 14).Value = -23745.31
$ws.Cells.Item(94, 8).Value = 35500
$ws.Cells.Item(94, 10).Value = 35500
$ws.Cells.Item(94, 12).Value = 35500
$ws.Cells.Item(94, 14).Value = -37302
$ws.Cells.Item(116, 8).Value = 18183278
$ws.Cells.Item(116, 9).Value = 30303796
$ws.Cells.Item(116, 11).Value = 30303796
$ws.Cells.Item(116, 13).Value = -30301502
$ws.Cells.Item(136, 9).Value = 4979.636
$ws.Cells.Item(136, 10).Value = 6939.125
$ws.Cells.Item(136, 11).Value = 14938.908
$ws.Cells.Item(136, 12).Value = 20817.375
$ws.Cells.Item(136, 13).Value = -12388.908
$ws.Cells.Item(136, 14).Value = -25917.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 18183278
$ws.Cells.Item(3, 9).Value = 30303796
$ws.Cells.Item(3, 11).Value = 30303796
$ws.Cells.Item(3, 13).Value = -30303682
$ws.Cells.Item(20, 8).Value = 3577.9375
$ws.Cells.Item(20, 9).Value = 2880
$ws.Cells.Item(20, 11).Value = 2880
$ws.Cells.Item(20, 13).Value = -2633
$ws.Cells.Item(94, 8).Value = 1190.303
$ws.Cells.Item(94, 9).Value = 902.8570999999999
$ws.Cells.Item(94, 11).Value = 902.8570999999999
$ws.Cells.Item(94, 13).Value = -451.8570999999999
$ws.Cells.Item(105, 8).Value = 52646090
$ws.Cells.Item(105, 9).Value = 58839576
$ws.Cells.Item(105, 11).Value = 58839576
$ws.Cells.Item(105, 13).Value = -58837829

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 37041024
$ws.Cells.Item(31, 9).Value = 100001930
$ws.Cells.Item(31, 10).Value = 5200.7646
$ws.Cells.Item(31, 11).Value = 100001930
$ws.Cells.Item(31, 12).Value = 5200.7646
$ws.Cells.Item(31, 13).Value = -100001635
$ws.Cells.Item(31, 14).Value = -5790.7646
$ws.Cells.Item(34, 8).Value = 37041024
$ws.Cells.Item(34, 9).Value = 100001930
$ws.Cells.Item(34, 10).Value = 5200.7646
$ws.Cells.Item(34, 11).Value = 100001930
$ws.Cells.Item(34, 12).Value = 5200.7646
$ws.Cells.Item(34, 13).Value = -100001728
$ws.Cells.Item(34, 14).Value = -5604.7646
$ws.Cells.Item(99, 8).Value = 3345.25
$ws.Cells.Item(99, 9).Value = 3441
$ws.Cells.Item(99, 10).Value = 3249.5
$ws.Cells.Item(99, 11).Value = 3441
$ws.Cells.Item(99, 12).Value = 3249.5
$ws.Cells.Item(99, 13).Value = -1943
$ws.Cells.Item(99, 14).Value = -6245.5
$ws.Cells.Item(107, 8).Value = 2052.4707
$ws.Cells.Item(107, 10).Value = 3367.4
$ws.Cells.Item(107, 12).Value = 3367.4
$ws.Cells.Item(107, 14).Value = -7207.4
$ws.Cells.Item(122, 8).Value = 112259.445
$ws.Cells.Item(122, 9).Value = 144048
$ws.Cells.Item(122, 11).Value = 432144
$ws.Cells.Item(122, 13).Value = -429694
$ws.Cells.Item(126, 8).Value = 3345.25
$ws.Cells.Item(126, 9).Value = 3441
$ws.Cells.Item(126, 10).Value = 3249.5
$ws.Cells.Item(126, 11).Value = 10323
$ws.Cells.Item(126, 12).Value = 9748.5
$ws.Cells.Item(126, 13).Value = -7853
$ws.Cells.Item(126, 14).Value = -14688.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 563.6667
$ws.Cells.Item(97, 9).Value = 492.57144
$ws.Cells.Item(97, 10).Value = 625.875
$ws.Cells.Item(97, 11).Value = 1477.71432
$ws.Cells.Item(97, 12).Value = 1877.625
$ws.Cells.Item(97, 13).Value = -981.71432
$ws.Cells.Item(97, 14).Value = -2869.625
$ws.Cells.Item(98, 8).Value = 199
$ws.Cells.Item(98, 10).Value = 199
$ws.Cells.Item(98, 12).Value = 597
$ws.Cells.Item(98, 14).Value = -3593
$ws.Cells.Item(107, 8).Value = 1755.0667
$ws.Cells.Item(107, 10).Value = 2317.6
$ws.Cells.Item(107, 12).Value = 6952.799999999999
$ws.Cells.Item(107, 14).Value = -10792.8
$ws.Cells.Item(113, 8).Value = 10000
$ws.Cells.Item(113, 10).Value = 10000
$ws.Cells.Item(113, 12).Value = 30000
$ws.Cells.Item(113, 14).Value = -34340
$ws.Cells.Item(132, 8).Value = 1779.8182
$ws.Cells.Item(132, 9).Value = 1323
$ws.Cells.Item(132, 11).Value = 11907
$ws.Cells.Item(132, 13).Value = -9377
$ws.Cells.Item(141, 8).Value = 5000
$ws.Cells.Item(141, 9).Value = 5000
$ws.Cells.Item(141, 11).Value = 15000
$ws.Cells.Item(141, 13).Value = -9820

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 20030
$ws.Cells.Item(52, 9).Value = 20030
$ws.Cells.Item(52, 11).Value = 20030
$ws.Cells.Item(52, 13).Value = -19771
$ws.Cells.Item(122, 8).Value = 5768.154
$ws.Cells.Item(122, 10).Value = 3969.2
$ws.Cells.Item(122, 12).Value = 11907.6
$ws.Cells.Item(122, 14).Value = -16807.6
$ws.Cells.Item(132, 8).Value = 4458.9062
$ws.Cells.Item(132, 9).Value = 2317.0625
$ws.Cells.Item(132, 11).Value = 6951.1875
$ws.Cells.Item(132, 13).Value = -4421.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(38, 8).Value = 25000
$ws.Cells.Item(38, 10).Value = 25000
$ws.Cells.Item(38, 12).Value = 25000
$ws.Cells.Item(38, 14).Value = -25820
$ws.Cells.Item(68, 8).Value = 2134.4546
$ws.Cells.Item(68, 9).Value = 1913.1666
$ws.Cells.Item(68, 11).Value = 1913.1666
$ws.Cells.Item(68, 13).Value = -1164.1666
$ws.Cells.Item(71, 8).Value = 2134.4546
$ws.Cells.Item(71, 9).Value = 1913.1666
$ws.Cells.Item(71, 11).Value = 9565.833000000001
$ws.Cells.Item(71, 13).Value = -5821.833000000001
$ws.Cells.Item(137, 8).Value = 99999
$ws.Cells.Item(137, 10).Value = 99999
$ws.Cells.Item(137, 12).Value = 99999
$ws.Cells.Item(137, 14).Value = -110199

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 97166.336
$ws.Cells.Item(46, 10).Value = 97166.336
$ws.Cells.Item(46, 12).Value = 97166.336
$ws.Cells.Item(46, 14).Value = -97628.336
$ws.Cells.Item(53, 8).Value = 19945.066
$ws.Cells.Item(53, 10).Value = 19945.066
$ws.Cells.Item(53, 12).Value = 19945.066
$ws.Cells.Item(53, 14).Value = -21159.066
$ws.Cells.Item(81, 8).Value = 6990.636
$ws.Cells.Item(81, 9).Value = 2442.5715
$ws.Cells.Item(81, 11).Value = 4885.143
$ws.Cells.Item(81, 13).Value = -3824.143
$ws.Cells.Item(84, 8).Value = 6990.636
$ws.Cells.Item(84, 9).Value = 2442.5715
$ws.Cells.Item(84, 11).Value = 24425.715
$ws.Cells.Item(84, 13).Value = -19121.715
$ws.Cells.Item(122, 8).Value = 5791.75
$ws.Cells.Item(122, 9).Value = 6047
$ws.Cells.Item(122, 11).Value = 18141
$ws.Cells.Item(122, 13).Value = -15691
$ws.Cells.Item(132, 8).Value = 4337.788
$ws.Cells.Item(132, 9).Value = 2445.2632
$ws.Cells.Item(132, 10).Value = 6906.2144
$ws.Cells.Item(132, 11).Value = 7335.7896
$ws.Cells.Item(132, 12).Value = 20718.6432
$ws.Cells.Item(132, 13).Value = -4805.7896
$ws.Cells.Item(132, 14).Value = -25778.6432
$ws.Cells.Item(134, 8).Value = 97166.336
$ws.Cells.Item(134, 10).Value = 97166.336
$ws.Cells.Item(134, 12).Value = 291499.008
$ws.Cells.Item(134, 14).Value = -296569.008
